$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header tweak: E1 "email" -> "Email" ---
$ws.Range("E1").Value = "Email"

# --- Row 2 : dari/hala -> clear email, error becomes "Invalid Email Error" ---
$ws.Range("E2").Value = ""
$ws.Range("G2").Value = "Invalid Email Error"

# --- Row 3 : overwritten with new person Fahim / Malika... ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Fahim"
$ws.Range("C3").Value = "Malikakkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkk"
$ws.Range("D3").Value = 1234567
$ws.Range("E3").Value = "Malika#gmail.com"
$ws.Range("F3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F3").Value = 36639
$ws.Range("G3").Value = "Invalid Email Error"

# --- Row 4 : overwritten with new person El hani / Hajar ---
$ws.Range("A4").Value = 123
$ws.Range("B4").Value = "El hani"
$ws.Range("C4").Value = "Hajar"
$ws.Range("D4").Value = 12345537
$ws.Range("E4").Value = "hajar@gmail.com"
$ws.Range("F4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F4").Value = 37588
$ws.Range("G4").Value = "Invalid Telephone Error"

# --- Row 5 : El Alami/Ahmed (unchanged data), error text updated ---
$ws.Range("G5").Value = "Invalid Telephone Error"

# --- Row 6 : overwritten with new person Alaoui / Fatima ---
$ws.Range("A6").Value = 21
$ws.Range("B6").Value = "Alaoui"
$ws.Range("C6").Value = "Fatima"
$ws.Range("D6").Value = 65434656
$ws.Range("E6").Value = "Fatima@gmail.com"
$ws.Range("F6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F6").Value = 37303
$ws.Range("G6").Value = "Invalid Telephone Error"

# --- Row 7 (new row) : gg / dfd (same data previously on row 3) ---
$ws.Range("A7").Value = 123
$ws.Range("B7").Value = "gg"
$ws.Range("C7").Value = "dfd"
$ws.Range("D7").Value = 1233
$ws.Range("E7").Value = "gg@gmail.com"
$ws.Range("F7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F7").Value = 40129
$ws.Range("G7").Value = "Invalid Telephone Error"

# --- Row 8 (new row) : hfg / dfdd (same data previously on row 4) ---
$ws.Range("A8").Value = 456
$ws.Range("B8").Value = "hfg"
$ws.Range("C8").Value = "dfdd"
$ws.Range("D8").Value = 12333222
$ws.Range("E8").Value = "dfdd@gmail.com"
$ws.Range("F8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("F8").Value = 37603
$ws.Range("G8").Value = "Invalid Telephone Error"

# --- Row 9 (new row) : El Alami / Ahmed (same data previously on row 5) ---
$ws.Range("A9").Value = 456
$ws.Range("B9").Value = "El Alami"
$ws.Range("C9").Value = "Ahmed"
$ws.Range("D9").Value = 61234567
$ws.Range("E9").Value = "Ahmed@gmail.com"
$ws.Range("F9").Value = "244/12/2003"
$ws.Range("G9").Value = "Invalid Date Error"

# --- Row 10 (new row) : Hamdaoui / Mohmed (same data previously on row 6) ---
$ws.Range("A10").Value = 789
$ws.Range("B10").Value = "Hamdaoui"
$ws.Range("C10").Value = "Mohmed"
$ws.Range("D10").Value = 123456789
$ws.Range("E10").Value = "mohmed@gmail.com"
$ws.Range("F10").Value = "12/38-2004"
$ws.Range("G10").Value = "Invalid Date Error"

# --- Row 11 (new row) : dari / hala duplicate (same data previously on row 7), email cleared ---
$ws.Range("A11").Value = 21
$ws.Range("B11").Value = "dari"
$ws.Range("C11").Value = "hala"
$ws.Range("D11").Value = 123456789
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "51/13/2001"
$ws.Range("G11").Value = "Invalid Date Error"
